$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 198.66667
$ws.Range("J2").Value = 410.25
$ws.Range("L2").Value = 410.25
$ws.Range("N2").Value = -636.25
# Row 5
$ws.Range("H5").Value = 199.77777
$ws.Range("I5").Value = 181
$ws.Range("J5").Value = 350
$ws.Range("K5").Value = 181
$ws.Range("L5").Value = 350
$ws.Range("M5").Value = -66
$ws.Range("N5").Value = -580
# Row 112
$ws.Range("H112").Value = 1997.5
$ws.Range("I112").Value = 1997.5
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 5992.5
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -4884.5
$ws.Range("N112").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 162.14285
$ws.Range("I4").Value = 115
$ws.Range("K4").Value = 115
$ws.Range("M4").Value = 1
# Row 5
$ws.Range("H5").Value = 75
$ws.Range("I5").Value = 100
$ws.Range("K5").Value = 100
$ws.Range("M5").Value = 12

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 75
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 15
# Row 29
$ws.Range("H29").Value = 690
$ws.Range("I29").Value = 637.5
$ws.Range("K29").Value = 637.5
$ws.Range("M29").Value = -348.5
# Row 57
$ws.Range("H57").Value = 95000
$ws.Range("J57").Value = 95000
$ws.Range("L57").Value = 95000
$ws.Range("N57").Value = -96440
# Row 75
$ws.Range("H75").Value = 38166.668
$ws.Range("I75").Value = 5000
$ws.Range("K75").Value = 5000
$ws.Range("M75").Value = -4064
# Row 78
$ws.Range("H78").Value = 38166.668
$ws.Range("I78").Value = 5000
$ws.Range("K78").Value = 15000
$ws.Range("M78").Value = -10320
# Row 94
$ws.Range("H94").Value = 606.8333
$ws.Range("I94").Value = 606.8333
$ws.Range("K94").Value = 606.8333
$ws.Range("M94").Value = -155.8333
# Row 136
$ws.Range("H136").Value = 95000
$ws.Range("J136").Value = 95000
$ws.Range("L136").Value = 95000
$ws.Range("N136").Value = -105200

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 86
$ws.Range("I7").Value = 129
$ws.Range("J7").Value = 55.285713
$ws.Range("K7").Value = 129
$ws.Range("L7").Value = 55.285713
$ws.Range("M7").Value = -16
$ws.Range("N7").Value = -281.285713
# Row 25
$ws.Range("H25").Value = 787.5
$ws.Range("J25").Value = 414.2857
$ws.Range("L25").Value = 414.2857
$ws.Range("N25").Value = -762.2857
# Row 58
$ws.Range("H58").Value = 3710.3333
$ws.Range("I58").Value = 1070
$ws.Range("J58").Value = 8991
$ws.Range("K58").Value = 1070
$ws.Range("L58").Value = 8991
$ws.Range("M58").Value = -867
$ws.Range("N58").Value = -9397
# Row 99
$ws.Range("H99").Value = 4115.3076
$ws.Range("I99").Value = 3826.25
$ws.Range("K99").Value = 3826.25
$ws.Range("M99").Value = -2328.25
# Row 105
$ws.Range("H105").Value = 1372.4
$ws.Range("I105").Value = 1372.4
$ws.Range("K105").Value = 1372.4
$ws.Range("M105").Value = 374.5999999999999
# Row 107
$ws.Range("H107").Value = 374.13333
$ws.Range("I107").Value = 352.72726
$ws.Range("K107").Value = 352.72726
$ws.Range("M107").Value = 1567.27274
# Row 126
$ws.Range("H126").Value = 4115.3076
$ws.Range("I126").Value = 3826.25
$ws.Range("K126").Value = 11478.75
$ws.Range("M126").Value = -9008.75
# Row 136
$ws.Range("H136").Value = 3710.3333
$ws.Range("I136").Value = 1070
$ws.Range("J136").Value = 8991
$ws.Range("K136").Value = 3210
$ws.Range("L136").Value = 26973
$ws.Range("M136").Value = -660
$ws.Range("N136").Value = -32073

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 36537756
$ws.Range("I4").Value = 56833710
$ws.Range("J4").Value = 5029.9
$ws.Range("K4").Value = 170501130
$ws.Range("L4").Value = 15089.7
$ws.Range("M4").Value = -170501018
$ws.Range("N4").Value = -15313.7
# Row 23
$ws.Range("H23").Value = 579.4545000000001
$ws.Range("J23").Value = 628.625
$ws.Range("L23").Value = 1885.875
$ws.Range("N23").Value = -2355.875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 66756.53
$ws.Range("I2").Value = 125041.625
$ws.Range("J2").Value = 145
$ws.Range("K2").Value = 125041.625
$ws.Range("L2").Value = 145
$ws.Range("M2").Value = -124928.625
$ws.Range("N2").Value = -371
# Row 3
$ws.Range("H3").Value = 13611848
$ws.Range("I3").Value = 19777858
$ws.Range("J3").Value = 3335166.5
$ws.Range("K3").Value = 19777858
$ws.Range("L3").Value = 3335166.5
$ws.Range("M3").Value = -19777742
$ws.Range("N3").Value = -3335398.5
# Row 5
$ws.Range("H5").Value = 3000
$ws.Range("I5").Value = 3000
$ws.Range("K5").Value = 3000
$ws.Range("M5").Value = -2888
# Row 7
$ws.Range("H7").Value = 12727809
$ws.Range("I7").Value = 17142858
$ws.Range("J7").Value = 5001475
$ws.Range("K7").Value = 17142858
$ws.Range("L7").Value = 5001475
$ws.Range("M7").Value = -17142746
$ws.Range("N7").Value = -5001699
# Row 8
$ws.Range("H8").Value = 12727809
$ws.Range("I8").Value = 17142858
$ws.Range("J8").Value = 5001475
$ws.Range("K8").Value = 17142858
$ws.Range("L8").Value = 5001475
$ws.Range("M8").Value = -17142719
$ws.Range("N8").Value = -5001753
# Row 11
$ws.Range("H11").Value = 12381238
$ws.Range("J11").Value = 6667333
$ws.Range("L11").Value = 6667333
$ws.Range("N11").Value = -6667611
# Row 14
$ws.Range("H14").Value = 11985.4
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 11985.4
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 11985.4
$ws.Range("N14").Value = -12321.4
$ws.Range("M14").ClearContents()
# Row 80
$ws.Range("H80").Value = 2008.5
$ws.Range("I80").Value = 1199.5
$ws.Range("J80").Value = 2817.5
$ws.Range("K80").Value = 1199.5
$ws.Range("L80").Value = 2817.5
$ws.Range("M80").Value = -201.5
$ws.Range("N80").Value = -4813.5
# Row 83
$ws.Range("H83").Value = 2008.5
$ws.Range("I83").Value = 1199.5
$ws.Range("J83").Value = 2817.5
$ws.Range("K83").Value = 5997.5
$ws.Range("L83").Value = 14087.5
$ws.Range("M83").Value = -1005.5
$ws.Range("N83").Value = -24071.5
# Row 97
$ws.Range("H97").Value = 2299.5
$ws.Range("I97").Value = 2299.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2299.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1803.5
$ws.Range("N97").ClearContents()
# Row 107
$ws.Range("H107").Value = 420.375
$ws.Range("I107").Value = 480
$ws.Range("K107").Value = 480
$ws.Range("M107").Value = 1440
# Row 122
$ws.Range("H122").Value = 716551.7
$ws.Range("I122").Value = 835727
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 2507181
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -2504731
$ws.Range("N122").Value = -9400

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7487.1113
$ws.Range("I7").Value = 6364.875
$ws.Range("K7").Value = 6364.875
$ws.Range("M7").Value = -6252.875
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
# Row 40
$ws.Range("H40").Value = 6491.143
$ws.Range("I40").Value = 3865
$ws.Range("J40").Value = 9992.666999999999
$ws.Range("K40").Value = 3865
$ws.Range("L40").Value = 9992.666999999999
$ws.Range("M40").Value = -3729
$ws.Range("N40").Value = -10264.667
# Row 68
$ws.Range("H68").Value = 9414.700000000001
$ws.Range("I68").Value = 8849.25
$ws.Range("J68").Value = 9791.666999999999
$ws.Range("K68").Value = 8849.25
$ws.Range("L68").Value = 9791.666999999999
$ws.Range("M68").Value = -8100.25
$ws.Range("N68").Value = -11289.667
# Row 71
$ws.Range("H71").Value = 9414.700000000001
$ws.Range("I71").Value = 8849.25
$ws.Range("J71").Value = 9791.666999999999
$ws.Range("K71").Value = 44246.25
$ws.Range("L71").Value = 48958.335
$ws.Range("M71").Value = -40502.25
$ws.Range("N71").Value = -56446.335
# Row 126
$ws.Range("H126").Value = 7487.1113
$ws.Range("I126").Value = 6364.875
$ws.Range("K126").Value = 19094.625
$ws.Range("M126").Value = -16624.625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 17
$ws.Range("H17").Value = 6502
$ws.Range("J17").Value = 5000
$ws.Range("L17").Value = 5000
$ws.Range("N17").Value = -5344
# Row 19
$ws.Range("H19").Value = 3000
$ws.Range("I19").Value = 3000
$ws.Range("J19").Value = 3000
$ws.Range("K19").Value = 3000
$ws.Range("L19").Value = 3000
$ws.Range("M19").Value = -2826
$ws.Range("N19").Value = -3348
# Row 126
$ws.Range("H126").Value = 5286.222
$ws.Range("I126").Value = 2017.75
$ws.Range("K126").Value = 6053.25
$ws.Range("M126").Value = -3583.25
